$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) updates remain plain text, matching the original formatting
$priceCells = @(2,3,4,5,6,7,8,9,12,13,16,17,18,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,40,41,42,43,44,45,47,48,50,51)
foreach ($r in $priceCells) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range('D2').Value = '28.007.31'
$ws.Range('E2').Value = '  +0.29%  '

$ws.Range('D3').Value = '1.859.53'
$ws.Range('E3').Value = '  -0.40%  '

$ws.Range('D4').Value = '1.004'
$ws.Range('E4').Value = '  +0.31%  '

$ws.Range('D5').Value = '311.81'
$ws.Range('E5').Value = '  -0.02%  '

$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  +0.23%  '

$ws.Range('D7').Value = '0.5099'
$ws.Range('E7').Value = '  +2.50%  '

$ws.Range('D8').Value = '0.3812'
$ws.Range('E8').Value = '  +0.26%  '

$ws.Range('D9').Value = '0.08265'
$ws.Range('E9').Value = '  -7.58%  '

$ws.Range('E10').Value = '  -0.61%  '

$ws.Range('E11').Value = '  +0.22%  '

$ws.Range('D12').Value = '6.199'
$ws.Range('E12').Value = '  -1.59%  '

$ws.Range('D13').Value = '1.864.49'
$ws.Range('E13').Value = '  +0.41%  '

$ws.Range('E14').Value = '  -0.82%  '

$ws.Range('E15').Value = '  -0.39%  '

$ws.Range('D16').Value = '1.004'
$ws.Range('E16').Value = '  +0.28%  '

$ws.Range('D17').Value = '0.00001095'
$ws.Range('E17').Value = '  -0.21%  '

$ws.Range('D18').Value = '90.46'
$ws.Range('E18').Value = '  -0.26%  '

$ws.Range('E19').Value = '  -0.02%  '

$ws.Range('E20').Value = '  -0.59%  '

$ws.Range('E21').Value = '  +0.11%  '

$ws.Range('D22').Value = '6.014'
$ws.Range('E22').Value = '  -0.99%  '

$ws.Range('D23').Value = '28.021.18'
$ws.Range('E23').Value = '  +0.28%  '

$ws.Range('D24').Value = '11.04'
$ws.Range('E24').Value = '  -2.97%  '

$ws.Range('D25').Value = '2.231'
$ws.Range('E25').Value = '  -2.33%  '

$ws.Range('D26').Value = '2.559'
$ws.Range('E26').Value = '  +1.62%  '

$ws.Range('D27').Value = '2.078.61'
$ws.Range('E27').Value = '  +0.27%  '

$ws.Range('D28').Value = '157.79'
$ws.Range('E28').Value = '  -0.02%  '

$ws.Range('D29').Value = '20.39'
$ws.Range('E29').Value = '  -1.31%  '

$ws.Range('D30').Value = '124.55'
$ws.Range('E30').Value = '  -1.15%  '

$ws.Range('D31').Value = '0.1060'
$ws.Range('E31').Value = '  +0.59%  '

$ws.Range('D32').Value = '1.035'
$ws.Range('E32').Value = '  -1.69%  '

$ws.Range('D33').Value = '5.596'
$ws.Range('E33').Value = '  +0.39%  '

$ws.Range('D34').Value = '3.599'
$ws.Range('E34').Value = '  +0.25%  '

$ws.Range('D35').Value = '9.645'
$ws.Range('E35').Value = '  +3.41%  '

$ws.Range('D36').Value = '0.06536'
$ws.Range('E36').Value = '  +0.09%  '

$ws.Range('D37').Value = '0.02414'
$ws.Range('E37').Value = '  +0.45%  '

$ws.Range('D38').Value = '0.2173'
$ws.Range('E38').Value = '  -0.31%  '

$ws.Range('E39').Value = '  +0.57%  '

$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '1.241'
$ws.Range('E40').Value = '  -2.66%  '

$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '0.6412'
$ws.Range('E41').Value = '  +0.91%  '

$ws.Range('D42').Value = '11.21'
$ws.Range('E42').Value = '  -3.61%  '

$ws.Range('D43').Value = '4.872'
$ws.Range('E43').Value = '  -0.35%  '

$ws.Range('D44').Value = '0.6097'
$ws.Range('E44').Value = '  +1.89%  '

$ws.Range('D45').Value = '13.03'
$ws.Range('E45').Value = '  -0.98%  '

$ws.Range('E46').Value = '  -0.48%  '

$ws.Range('D47').Value = '3.653'
$ws.Range('E47').Value = '  -0.37%  '

$ws.Range('D48').Value = '1.980'
$ws.Range('E48').Value = '  +0.91%  '

$ws.Range('E49').Value = '  -1.17%  '

$ws.Range('D50').Value = '120.77'

$ws.Range('D51').Value = '79.09'
$ws.Range('E51').Value = '  +1.51%  '
